# DVT Led HW Change
$wb = $excel.ActiveWorkbook
$wsCmd = $wb.Worksheets.Item("commands_list")
$wsLed = $wb.Worksheets.Item("leds_names")

# --- 1. Remove the obsolete "blue"-variant rows on leds_names ---
# Rows (before edit) containing: ACTION_b(1), Lantern_b(7), PLAYPAUSE_b(10),
# USB_b(13), USB_g(14), VOL+_b(16), VOL-_b(19).
# (CONNECT_b on row 4 is kept.) Delete bottom-up so earlier row numbers stay valid.
$rowsToDelete = 19,16,14,13,10,7,1
foreach ($r in $rowsToDelete) {
    $wsLed.Rows.Item($r).Delete()
}

# --- 2. Re-apply the group header labels (lost with their deleted cell) to the new top row
#         of each shrunk merged group. ("Connect Key" on row 3 survived the delete untouched.) ---
$wsLed.Range("A1").Value = "Action Key"
$wsLed.Range("A6").Value = "Lantern Key"
$wsLed.Range("A8").Value = "Pause Key"
$wsLed.Range("A10").Value = "USB"
$wsLed.Range("A11").Value = "Volume+ Key"
$wsLed.Range("A13").Value = "Volume- Key"

# --- 3. Column C on leds_names: widen / autofit ---
$wsLed.Columns.Item(3).ColumnWidth = 39.375
$wsLed.Columns.Item(3).AutoFit()

# --- 4. View / selection state ---
$wsCmd.Activate()
$wsCmd.Application.ActiveWindow.ScrollRow = 31
$wsCmd.Range("F20").Select()

$wsLed.Activate()
$wsLed.Application.ActiveWindow.ScrollRow = 1
$wsLed.Range("B1").Select()
